$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# This workbook tracks handback status for two files that were handed off:
#   96bed73a-da10-4345-97b5-71c71852a50e.md   (kept)
#   f0dc634e-0397-4e76-8a16-684b965cd40a.md   (now fully handed back -> removed
#                                              from the report, and the
#                                              surviving row's handoff/handback
#                                              timestamps are refreshed)
#
# For every sheet this means: drop row 3 (and its hyperlinks), and shrink the
# used range accordingly. On the two per-locale sheets, row 2 also gets new
# "Correspond Handoff/Handback Datetime" values reflecting the newer run.
# ---------------------------------------------------------------------------

# ---- Sheet "Overview" ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Hyperlinks.Delete()
$wsOverview.Rows.Item(3).Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/e140c29a671146e679931ae437764f38d830fa77/e2e/96bed73a-da10-4345-97b5-71c71852a50e.md", "", "", "96bed73a-da10-4345-97b5-71c71852a50e.md")

# ---- Sheet "zh-cn" ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-24 02:55:13"
$wsZhCn.Range("H2").Value = "2016-03-24 02:56:04"
$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Rows.Item(3).Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/e140c29a671146e679931ae437764f38d830fa77/e2e/96bed73a-da10-4345-97b5-71c71852a50e.md", "", "", "96bed73a-da10-4345-97b5-71c71852a50e.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7a9b030926ac1b625269e73ad0c081302c3cfca6/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/96bed73a-da10-4345-97b5-71c71852a50e.346d42a83ba9f8cb2cc7bf27941f807f5f940db2.zh-cn.xlf", "", "", "96bed73a-da10-4345-97b5-71c71852a50e.346d42a83ba9f8cb2cc7bf27941f807f5f940db2.zh-cn.xlf")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/844a60f95d768e8e99e0dfe75fc9694b9df2e5e9/e2e/96bed73a-da10-4345-97b5-71c71852a50e.md", "", "", "96bed73a-da10-4345-97b5-71c71852a50e.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/539faadb096e7245488654e02857cd4af7b1a335/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/96bed73a-da10-4345-97b5-71c71852a50e.346d42a83ba9f8cb2cc7bf27941f807f5f940db2.zh-cn.xlf", "", "", "96bed73a-da10-4345-97b5-71c71852a50e.346d42a83ba9f8cb2cc7bf27941f807f5f940db2.zh-cn.xlf")

# ---- Sheet "de-de" ----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-24 02:55:17"
$wsDeDe.Range("H2").Value = "2016-03-24 02:56:11"
$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Rows.Item(3).Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/e140c29a671146e679931ae437764f38d830fa77/e2e/96bed73a-da10-4345-97b5-71c71852a50e.md", "", "", "96bed73a-da10-4345-97b5-71c71852a50e.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/41bdd46b98169a248611291899cf0c9bb78625fb/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/96bed73a-da10-4345-97b5-71c71852a50e.346d42a83ba9f8cb2cc7bf27941f807f5f940db2.de-de.xlf", "", "", "96bed73a-da10-4345-97b5-71c71852a50e.346d42a83ba9f8cb2cc7bf27941f807f5f940db2.de-de.xlf")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/2fd1e238981044a33036a51cf1c64cd4e528352e/e2e/96bed73a-da10-4345-97b5-71c71852a50e.md", "", "", "96bed73a-da10-4345-97b5-71c71852a50e.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/b585f1dca2112b1e42099da51b39bb6db7d14913/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/96bed73a-da10-4345-97b5-71c71852a50e.346d42a83ba9f8cb2cc7bf27941f807f5f940db2.de-de.xlf", "", "", "96bed73a-da10-4345-97b5-71c71852a50e.346d42a83ba9f8cb2cc7bf27941f807f5f940db2.de-de.xlf")
